$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 189, pushing the existing
# rows 189:295 down to 191:297.
$ws.Range("A189:A190").EntireRow.Insert()

# Populate the first new row (189) - "Primera" quality record
$ws.Range("A189").Value = 6
$ws.Range("B189").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C189").Value = "Metropolitana"
$ws.Range("D189").Value = 44529
$ws.Range("E189").Value = 13
$ws.Range("F189").Value = 100112052
$ws.Range("G189").Value = "Albahaca"
$ws.Range("H189").Value = "Sin especificar"
$ws.Range("I189").Value = "Primera"
$ws.Range("J189").Value = 130
$ws.Range("K189").Value = 6000
$ws.Range("L189").Value = 6000
$ws.Range("M189").Value = 6000
$ws.Range("N189").Value = "`$/docena de matas"
$ws.Range("O189").Value = "Región Metropolitana"
$ws.Range("P189").Value = 1000
$ws.Range("Q189").Value = 6
$ws.Range("R189").Value = "Hortaliza"

# Populate the second new row (190) - "Segunda" quality record
$ws.Range("A190").Value = 6
$ws.Range("B190").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C190").Value = "Metropolitana"
$ws.Range("D190").Value = 44529
$ws.Range("E190").Value = 13
$ws.Range("F190").Value = 100112052
$ws.Range("G190").Value = "Albahaca"
$ws.Range("H190").Value = "Sin especificar"
$ws.Range("I190").Value = "Segunda"
$ws.Range("J190").Value = 60
$ws.Range("K190").Value = 5000
$ws.Range("L190").Value = 5000
$ws.Range("M190").Value = 5000
$ws.Range("N190").Value = "`$/docena de matas"
$ws.Range("O190").Value = "Región Metropolitana"
$ws.Range("P190").Value = 833
$ws.Range("Q190").Value = 6
$ws.Range("R190").Value = "Hortaliza"
